# Add three new header columns (D, E, F) to Sheet1 describing additional
# ORG_ROOM fields: ORG_ROOM_IDENOLD, ORG_ROOM_IDENNEW, ORG_ROOM_STATUS.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "ORG_ROOM_IDENOLD"
$ws.Range("E1").Value = "ORG_ROOM_IDENNEW"
$ws.Range("F1").Value = "ORG_ROOM_STATUS"

# E1/F1 pick up the same left-aligned style used by the existing header
# cells (A1:C1); D1 keeps the default (unstyled) formatting, matching the
# way the new data was pasted in on the original upload.
$xlLeft = -4131
$ws.Range("E1").HorizontalAlignment = $xlLeft
$ws.Range("F1").HorizontalAlignment = $xlLeft

[void]$ws.Range("F2").Select()
